# Extend the per-person "group" columns from column DL out to column EH.
# The sheet already repeats each row's group label (taken from column C)
# across columns C..DL; this mirrors that same label into the 22 new
# columns DM..EH for every data row (2-15). Row 10 has no group label
# (its C:DL cells are blank), so the newly added cells there stay blank
# too - we still materialize them so the used range / dimension extends
# exactly like the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$groupByRow = @{
    2  = "group3"
    3  = "group1"
    4  = "group3"
    5  = "group3"
    6  = "group1"
    7  = "group1"
    8  = "group2"
    9  = "group2"
    10 = ""
    11 = "group2"
    12 = "group2"
    13 = "group2"
    14 = "group2"
    15 = "group1"
}

$newCols = @("DM","DN","DO","DP","DQ","DR","DS","DT","DU","DV","DW","DX","DY","DZ", `
             "EA","EB","EC","ED","EE","EF","EG","EH")

foreach ($row in 2..15) {
    $value = $groupByRow[$row]
    foreach ($col in $newCols) {
        $cell = $ws.Range("$col$row")
        if ($value -eq "") {
            # Force the otherwise-empty cell to materialize (matching the
            # blank placeholder cells already present at C10:DL10) since a
            # plain empty-string value write is treated as a no-op. Touching
            # the font (re-asserting the sheet's own Arial font) is a
            # harmless no-op formatting-wise but makes the cell "real".
            $cell.Font.Name = "Arial"
        } else {
            $cell.Value = $value
        }
    }
}
